# miRNA best correl/worst correl gene tables
#
# - Fix the mRNA sheet's tab name (drop the stray trailing parenthesis).
# - Move the selection on the miRNA sheet from O7 to J13.
# - Keep the mRNA sheet ("Comparison Tables mRNA") the active/selected tab,
#   matching the original tabSelected state.

$wb = $excel.ActiveWorkbook

$wsMrna  = $wb.Worksheets.Item(1)
$wsMirna = $wb.Worksheets.Item(2)

# Rename "Comparison Tables mRNA)" -> "Comparison Tables mRNA"
$wsMrna.Name = "Comparison Tables mRNA"

# Update the selected/active cell on the miRNA sheet to J13.
[void]$wsMirna.Activate()
[void]$wsMirna.Range("J13").Select()

# Restore the mRNA sheet as the active tab (it was, and should remain,
# tabSelected="1").
[void]$wsMrna.Activate()
